# complex_validate_test.xlsx edit
#
# Commit message: "Change 'assign' to use 'calculation' as the column name
# for the formula to assign to the name. Remove 'default' from prompts.
# Can specify an value default on the model element, or use 'assign'."
#
# In sheet terms: the "survey"/"section1"/"section2" sheets each had a
# header column named 'default' (column R, right after 'choice_filter').
# That column is removed entirely on all three sheets, which shifts the
# remaining trailing column(s) ('validation_tags' / 'hideInContents') one
# slot to the left.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("survey")
$ws2 = $wb.Worksheets.Item("section1")
$ws3 = $wb.Worksheets.Item("section2")

# Drop the 'default' column on each of the three sheets. Excel shifts
# everything to its right one column to the left automatically, which is
# exactly the R->(old S/T) shift seen in the diff.
$ws1.Columns("R").Delete()
$ws2.Columns("R").Delete()
$ws3.Columns("R").Delete()

# Match the resulting selection/view state on each sheet.
$ws1.Columns("R").Select()
$ws2.Columns("R").Select()
$ws3.Range("P10").Select()

# "section2" ends up the active/selected tab.
$ws3.Activate()
